$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.253.57"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.408.68"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'232.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "'618.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("D7").Value = "'1.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.15%  "
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'0.980"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.73%  "
$ws.Range("D11").Value = "3.405.34"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").Value = "'43.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.60%  "
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("E14").Value = "  +4.93%  "
$ws.Range("D15").Value = "93.086.42"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "4.058.74"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "'8.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.63%  "
$ws.Range("D19").Value = "3.422.92"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").Value = "'17.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.08%  "
$ws.Range("D21").Value = "'11.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.32%  "
$ws.Range("D22").Value = "'0.503"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.41%  "
$ws.Range("E23").Value = "  +8.38%  "
$ws.Range("D24").Value = "'496.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "'6.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.53%  "
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("D27").Value = "'90.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("E28").Value = "  +5.13%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'11.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +5.74%  "
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").Value = "'0.175"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("D35").Value = "'0.548"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("D36").Value = "'28.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").Value = "'557.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.59%  "
$ws.Range("D38").Value = "'7.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").Value = "'0.149"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("D42").Value = "'0.894"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").Value = "'1.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "'3.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").Value = "'0.0413"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.82%  "
$ws.Range("D47").Value = "'5.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "'52.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("E51").Value = "  -1.32%  "
